# Append 11 new NBA game rows (rows 848-858) to Sheet1, mirroring the
# existing table layout:
#   A: Away team   B: Away Pts   C: Home team   D: Home Pts
#   E: Overtime    F: Attend.    G: Arena
#   H: Win         I: Loss

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Milwaukee Bucks",      119, "Philadelphia 76ers",    98, "No", 17832, "Wells Fargo Center",    "Milwaukee Bucks",      "Philadelphia 76ers"),
    @("Los Angeles Lakers",   113, "Phoenix Suns",          123, "No", 17832, "Footprint Center",      "Phoenix Suns",         "Los Angeles Lakers"),
    @("Dallas Mavericks",     111, "Indiana Pacers",        133, "No", 17832, "Gainbridge Fieldhouse", "Indiana Pacers",       "Dallas Mavericks"),
    @("Cleveland Cavaliers",  114, "Washington Wizards",    105, "No", 17832, "Capital One Arena",     "Cleveland Cavaliers",  "Washington Wizards"),
    @("Orlando Magic",         92, "Atlanta Hawks",         109, "No", 17832, "State Farm Arena",      "Atlanta Hawks",        "Orlando Magic"),
    @("Denver Nuggets",       119, "Golden State Warriors", 103, "No", 17832, "Chase Center",          "Denver Nuggets",       "Golden State Warriors"),
    @("Oklahoma City Thunder",123, "Houston Rockets",       110, "No", 17832, "Toyota Center",         "Oklahoma City Thunder","Houston Rockets"),
    @("Chicago Bulls",        114, "New Orleans Pelicans",  106, "No", 17832, "Smoothie King Center",  "Chicago Bulls",        "New Orleans Pelicans"),
    @("San Antonio Spurs",    109, "Utah Jazz",              128,"No", 17832, "Delta Center",          "Utah Jazz",            "San Antonio Spurs"),
    @("Charlotte Hornets",     93, "Portland Trail Blazers", 80,"No", 17832, "Moda Center",           "Charlotte Hornets",    "Portland Trail Blazers"),
    @("Sacramento Kings",     123, "Los Angeles Clippers",  107, "No", 17832, "Crypto.com Arena",      "Sacramento Kings",     "Los Angeles Clippers")
)

$startRow = 848
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}

# Scroll/select to mirror the author's final view position (row 829 at the
# top, new blank row A859 selected as the next entry point).
$excel.ActiveWindow.ScrollRow = 829
$ws.Range("A859").Select()
